# Update the two-digit division worksheet numbers.
#
# Each entry is (row, col, newText) referring to cells in the single table
# in the document. Rows 1, 5, 9, 13, 17 (1-based) hold the "N÷N=" problems;
# the rows in between are blank answer rows.
#
# NOTE: several "old" values are duplicated across cells (e.g. "41÷3="
# appears twice) and several "new" values collide with "old" values located
# in other cells (e.g. one cell's old value "40÷8=" becomes "56÷8=", while
# another cell's old value is itself "56÷8="). A plain document-wide
# Find/Replace (even when started from a Range scoped to a single cell)
# is not reliably confined to that range in this host, so instead we
# directly overwrite each cell's text range by precise character offsets,
# which keeps every edit isolated to its own cell and preserves the
# existing run formatting.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cellRange = $table.Cell($row, $col).Range
    # Cell range Text includes the trailing cell-mark character(s); trim
    # the last one off so we only overwrite the visible "N÷N=" content.
    $r = $d.Range($cellRange.Start, $cellRange.End - 1)
    $r.Text = $newText
}

$changes = @(
    @(1, 1, "31÷9="),
    @(1, 2, "82÷4="),
    @(1, 3, "19÷9="),
    @(1, 4, "56÷8="),
    @(1, 5, "93÷9="),

    @(5, 1, "64÷5="),
    @(5, 2, "75÷6="),
    @(5, 3, "40÷5="),
    @(5, 4, "23÷6="),
    @(5, 5, "49÷3="),

    @(9, 1, "35÷3="),
    @(9, 2, "42÷5="),
    @(9, 3, "28÷9="),
    @(9, 4, "62÷2="),
    @(9, 5, "90÷6="),

    @(13, 1, "92÷6="),
    @(13, 2, "42÷4="),
    @(13, 3, "28÷7="),
    @(13, 4, "81÷4="),
    @(13, 5, "44÷7="),

    @(17, 1, "93÷6="),
    @(17, 2, "38÷6="),
    @(17, 3, "85÷8="),
    @(17, 4, "22÷2="),
    @(17, 5, "25÷9=")
)

foreach ($change in $changes) {
    $row = $change[0]
    $col = $change[1]
    $new = $change[2]

    Set-CellText $t $row $col $new
}
